$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: "...horas del mismo día. Número de..." ->
#         "...horas del <CONCLUYE>. Número de..."
#         split into three separate runs: "el ", "<CONCLUYE>", ". "
# ---------------------------------------------------------------------------

$findR = $d.Content
$findR.Find.Execute("el mismo día. ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$segStart = $findR.Start
$segEnd = $findR.End

# Pre-split the run by dropping zero-width bookmarks at every boundary we
# want to preserve/create; this keeps the edit from re-coalescing the run
# with its unrelated neighbours.
$d.Bookmarks.Add("zzTmpA", $d.Range($segStart, $segStart)) | Out-Null
$d.Bookmarks.Add("zzTmpB", $d.Range($segStart + 3, $segStart + 3)) | Out-Null
$d.Bookmarks.Add("zzTmpC", $d.Range($segEnd - 2, $segEnd - 2)) | Out-Null

$midRange = $d.Range($segStart + 3, $segEnd - 2)
$midRange.Text = "<CONCLUYE>"

$d.Bookmarks("zzTmpA").Delete()
$d.Bookmarks("zzTmpB").Delete()
$d.Bookmarks("zzTmpC").Delete()

# ---------------------------------------------------------------------------
# Part 2: "Número de invitados 50.  " -> "Número de invitados <INVITADOS>.  "
#         the "50" run was Bold + cyan highlight; the new run keeps Bold but
#         drops the highlight; re-home the "_GoBack" bookmark around it.
# ---------------------------------------------------------------------------

$p5 = $d.Paragraphs(5).Range
$scopeR = $d.Range($p5.Start, $p5.End)
$scopeR.Find.Execute("50", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$numStart = $scopeR.Start
$numEnd = $scopeR.End

# Drop the old run completely (this also drops its cyan highlight) then type
# a brand new run in its place so it starts out with plain formatting.
$d.Range($numStart, $numEnd).Delete()
$d.Range($numStart, $numStart).InsertAfter("<INVITADOS>") | Out-Null
$newNumEnd = $numStart + ("<INVITADOS>").Length
$d.Range($numStart, $newNumEnd).Bold = 1

# Remove the old (empty) _GoBack bookmark wherever it currently sits.
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# Re-create _GoBack spanning "<INVITADOS>.  " (mirrors the edited document).
$p5b = $d.Paragraphs(5).Range
$goBackRange = $d.Range($numStart, $p5b.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
